$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H26").Value = 10000
$ws.Range("J26").Value = 10000
$ws.Range("L26").Value = 10000
$ws.Range("N26").Value = -10688
$ws.Range("H41").Value = 875.25
$ws.Range("I41").Value = 1000
$ws.Range("K41").Value = 1000
$ws.Range("M41").Value = -560
$ws.Range("H125").Value = 11365592
$ws.Range("I125").Value = 41667830
$ws.Range("J125").Value = 2251.75
$ws.Range("K125").Value = 375010470
$ws.Range("L125").Value = 20265.75
$ws.Range("M125").Value = -375008010
$ws.Range("N125").Value = -25185.75
$ws.Range("H132").Value = 2090.862
$ws.Range("I132").Value = 2134.532
$ws.Range("K132").Value = 6403.596
$ws.Range("M132").Value = -3873.596
$ws.Range("H135").Value = 29413498
$ws.Range("I135").Value = 33334632
$ws.Range("K135").Value = 300011688
$ws.Range("M135").Value = -300009153
$ws.Range("H137").Value = 2875146.2
$ws.Range("I137").Value = 6945496.5
$ws.Range("J137").Value = 1957.9412
$ws.Range("K137").Value = 20836489.5
$ws.Range("L137").Value = 5873.8236
$ws.Range("M137").Value = -20833939.5
$ws.Range("N137").Value = -10973.8236
$ws.Range("H138").Value = 2896.02
$ws.Range("I138").Value = 1097.6364
$ws.Range("J138").Value = 3403.2563
$ws.Range("K138").Value = 3292.9092
$ws.Range("L138").Value = 10209.7689
$ws.Range("M138").Value = 1847.0908
$ws.Range("N138").Value = -20489.7689
$ws.Range("H140").Value = 75836
$ws.Range("J140").Value = 75836
$ws.Range("L140").Value = 75836
$ws.Range("N140").Value = -86196
$ws.Range("H141").Value = 3221.32
$ws.Range("I141").Value = 1482.4
$ws.Range("J141").Value = 10177
$ws.Range("K141").Value = 4447.200000000001
$ws.Range("L141").Value = 30531
$ws.Range("M141").Value = 732.7999999999993
$ws.Range("N141").Value = -40891
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10428797
$ws.Range("I32").Value = 11122371
$ws.Range("J32").Value = 25185
$ws.Range("K32").Value = 11122371
$ws.Range("L32").Value = 25185
$ws.Range("M32").Value = -11122084
$ws.Range("N32").Value = -25759
$ws.Range("H74").Value = 11366042
$ws.Range("I74").Value = 1777.8846
$ws.Range("J74").Value = 27781090
$ws.Range("K74").Value = 1777.8846
$ws.Range("L74").Value = 27781090
$ws.Range("M74").Value = -903.8846000000001
$ws.Range("N74").Value = -27782838
$ws.Range("H77").Value = 11366042
$ws.Range("I77").Value = 1777.8846
$ws.Range("J77").Value = 27781090
$ws.Range("K77").Value = 8889.423000000001
$ws.Range("L77").Value = 138905450
$ws.Range("M77").Value = -4521.423000000001
$ws.Range("N77").Value = -138914186
$ws.Range("H102").Value = 2670
$ws.Range("I102").Value = 2670
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 2670
$ws.Range("L102").Value = 0
$ws.Range("M102").ClearContents()
$ws.Range("N102").Value = -1048
$ws.Range("H132").Value = 1712942.4
$ws.Range("I132").Value = 2844.6897
$ws.Range("K132").Value = 8534.069100000001
$ws.Range("M132").Value = -6004.069100000001
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 27028738
$ws.Range("I20").Value = 1596.238
$ws.Range("J20").Value = 62501864
$ws.Range("K20").Value = 1596.238
$ws.Range("L20").Value = 62501864
$ws.Range("M20").Value = -1349.238
$ws.Range("N20").Value = -62502358
$ws.Range("H134").Value = 3165.0312
$ws.Range("I134").Value = 3195.5
$ws.Range("J134").Value = 3098
$ws.Range("K134").Value = 9586.5
$ws.Range("L134").Value = 9294
$ws.Range("M134").Value = -7051.5
$ws.Range("N134").Value = -14364
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5265.3076
$ws.Range("I31").Value = 1782.4286
$ws.Range("J31").Value = 6548.4736
$ws.Range("K31").Value = 1782.4286
$ws.Range("L31").Value = 6548.4736
$ws.Range("M31").Value = -1487.4286
$ws.Range("N31").Value = -7138.4736
$ws.Range("H34").Value = 5265.3076
$ws.Range("I34").Value = 1782.4286
$ws.Range("J34").Value = 6548.4736
$ws.Range("K34").Value = 1782.4286
$ws.Range("L34").Value = 6548.4736
$ws.Range("M34").Value = -1580.4286
$ws.Range("N34").Value = -6952.4736
$ws.Range("H58").Value = 1261.7858
$ws.Range("I58").Value = 1005.1
$ws.Range("J58").Value = 1903.5
$ws.Range("K58").Value = 1005.1
$ws.Range("L58").Value = 1903.5
$ws.Range("M58").Value = -802.1
$ws.Range("N58").Value = -2309.5
$ws.Range("H132").Value = 21507526
$ws.Range("I132").Value = 31252008
$ws.Range("J132").Value = 11113412
$ws.Range("K132").Value = 93756024
$ws.Range("L132").Value = 33340236
$ws.Range("M132").Value = -93753494
$ws.Range("N132").Value = -33345296
$ws.Range("H136").Value = 1261.7858
$ws.Range("I136").Value = 1005.1
$ws.Range("J136").Value = 1903.5
$ws.Range("K136").Value = 3015.3
$ws.Range("L136").Value = 5710.5
$ws.Range("M136").Value = -465.3000000000002
$ws.Range("N136").Value = -10810.5
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 144.73334
$ws.Range("J2").Value = 171
$ws.Range("L2").Value = 1026
$ws.Range("N2").Value = -1252
$ws.Range("H5").Value = 1004.8333
$ws.Range("I5").Value = 737.5909
$ws.Range("K5").Value = 2212.7727
$ws.Range("M5").Value = -2100.7727
$ws.Range("H12").Value = 65.791664
$ws.Range("I12").Value = 35.3125
$ws.Range("K12").Value = 105.9375
$ws.Range("M12").Value = 67.0625
$ws.Range("H93").Value = 4981.25
$ws.Range("I93").Value = 4922
$ws.Range("K93").Value = 14766
$ws.Range("M93").Value = -12894
$ws.Range("H108").Value = 1284
$ws.Range("I108").Value = 940.8
$ws.Range("K108").Value = 2822.4
$ws.Range("M108").Value = 57.60000000000036
$ws.Range("H135").Value = 1004.8333
$ws.Range("I135").Value = 737.5909
$ws.Range("K135").Value = 6638.3181
$ws.Range("M135").Value = -4103.3181
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 5000
$ws.Range("I5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("M5").ClearContents()
$ws.Range("H70").Value = 5463.5
$ws.Range("I70").Value = 5413.1333
$ws.Range("J70").Value = 5571.4287
$ws.Range("K70").Value = 5413.1333
$ws.Range("L70").Value = 5571.4287
$ws.Range("M70").Value = -5143.1333
$ws.Range("N70").Value = -6111.4287
$ws.Range("H73").Value = 5463.5
$ws.Range("I73").Value = 5413.1333
$ws.Range("J73").Value = 5571.4287
$ws.Range("K73").Value = 5413.1333
$ws.Range("L73").Value = 5571.4287
$ws.Range("M73").Value = -4477.1333
$ws.Range("N73").Value = -7443.4287
$ws.Range("H132").Value = 30308960
$ws.Range("I132").Value = 58831812
$ws.Range("J132").Value = 3428.625
$ws.Range("K132").Value = 176495436
$ws.Range("L132").Value = 10285.875
$ws.Range("M132").Value = -176492906
$ws.Range("N132").Value = -15345.875
$ws.Range("H135").Value = 55000
$ws.Range("J135").Value = 55000
$ws.Range("L135").Value = 55000
$ws.Range("N135").Value = -65140
$ws.Range("H136").Value = 28786.77
$ws.Range("J136").Value = 25929.818
$ws.Range("L136").Value = 77789.454
$ws.Range("N136").Value = -82889.454
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4142
$ws.Range("I40").Value = 3652.5
$ws.Range("J40").Value = 4876.25
$ws.Range("K40").Value = 3652.5
$ws.Range("L40").Value = 4876.25
$ws.Range("M40").Value = -3516.5
$ws.Range("N40").Value = -5148.25
$ws.Range("H82").Value = 2929.4
$ws.Range("I82").Value = 3291.6667
$ws.Range("J82").Value = 2386
$ws.Range("K82").Value = 3291.6667
$ws.Range("L82").Value = 2386
$ws.Range("M82").Value = -2930.6667
$ws.Range("N82").Value = -3108
$ws.Range("H85").Value = 2929.4
$ws.Range("I85").Value = 3291.6667
$ws.Range("J85").Value = 2386
$ws.Range("K85").Value = 3291.6667
$ws.Range("L85").Value = 2386
$ws.Range("M85").Value = -2043.6667
$ws.Range("N85").Value = -4882
$ws.Range("H132").Value = 3375.1667
$ws.Range("I132").Value = 3000.1667
$ws.Range("J132").Value = 3750.1667
$ws.Range("K132").Value = 9000.500100000001
$ws.Range("L132").Value = 11250.5001
$ws.Range("M132").Value = -6470.500100000001
$ws.Range("N132").Value = -16310.5001
$ws.Range("H140").Value = 61358.5
$ws.Range("J140").Value = 61358.5
$ws.Range("L140").Value = 61358.5
$ws.Range("N140").Value = -71718.5
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H94").Value = 64165
$ws.Range("J94").Value = 64165
$ws.Range("L94").Value = 64165
$ws.Range("N94").Value = -65967
$ws.Range("H132").Value = 6483870.5
$ws.Range("I132").Value = 2464.6072
$ws.Range("J132").Value = 17159128
$ws.Range("K132").Value = 7393.821599999999
$ws.Range("L132").Value = 51477384
$ws.Range("M132").Value = -4863.821599999999
$ws.Range("N132").Value = -51482444
$ws.Range("H136").Value = 4450.206
$ws.Range("I136").Value = 4171.5
$ws.Range("J136").Value = 4763.75
$ws.Range("K136").Value = 12514.5
$ws.Range("L136").Value = 14291.25
$ws.Range("M136").Value = -9964.5
$ws.Range("N136").Value = -19391.25
